$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 730.6461864673427
$ws.Range("C2").Value = 2905.942428134453
$ws.Range("D2").Value = 2192.830618823093
